$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 41 - this shifts existing rows 41..144 down to 42..145
$ws.Rows.Item(41).Insert()

# Populate the newly inserted row 41 with the new data record
$ws.Cells.Item(41, 1).Value = 8
$ws.Cells.Item(41, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(41, 3).Value = "Coquimbo"
$ws.Cells.Item(41, 4).Value = 44914
$ws.Cells.Item(41, 5).Value = 4
$ws.Cells.Item(41, 6).Value = 100112052
$ws.Cells.Item(41, 7).Value = "Albahaca"
$ws.Cells.Item(41, 8).Value = "Sin especificar"
$ws.Cells.Item(41, 9).Value = "Primera"
$ws.Cells.Item(41, 10).Value = 1100
$ws.Cells.Item(41, 11).Value = 3000
$ws.Cells.Item(41, 12).Value = 4000
$ws.Cells.Item(41, 13).Value = 3500
$ws.Cells.Item(41, 14).Value = "$/paquete"
$ws.Cells.Item(41, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(41, 16).Value = 3500
$ws.Cells.Item(41, 17).Value = 1
$ws.Cells.Item(41, 18).Value = "Hortaliza"
